$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "TC001"
$ws.Range("B3").Value = "Creating the TEAM Workspace"
$ws.Range("C3").Value = "TEAM Workspace Should be created successfully and approved"
$ws.Range("D3").Value = "TEAM workspace is created successfully and approved"
$ws.Range("E3").Value = "Pass"
